# Upload_file_to_GitHub.docx - "First commit to the branch"
#
# The document ends with four paragraphs right before the final <w:sectPr>:
#   ... "               ?????????? (master)"
#   [A] empty Plain Text paragraph
#   [B] empty Plain Text paragraph
#   [C] paragraph holding only the _GoBack bookmark (no pPr -> Normal style)
#
# The edit:
#   - [A] gains the run "$ git status"
#   - [B] gains the run "$ git add ." and the _GoBack bookmark (moved from [C])
#   - two brand new empty Plain Text paragraphs are appended after [B]
#   - [C] survives as a final, completely empty paragraph (no bookmark anymore)

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
$paraA = $d.Paragraphs.Item($count - 2)   # empty Plain Text paragraph -> "$ git status"
$paraB = $d.Paragraphs.Item($count - 1)   # empty Plain Text paragraph -> "$ git add ."
$paraC = $d.Paragraphs.Item($count)       # bookmark-only paragraph (_GoBack)

$pkgNs = 'xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"'
$wNs   = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-PkgXml($bodyInner) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package ' + $pkgNs + `
        '><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $wNs + `
        '><w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

$rFonts = '<w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>'

# Remove the existing _GoBack bookmark; it will be re-created inside paragraph B.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# Paragraph A: insert the "$ git status" run without touching its own pPr/mark.
$rngA = $d.Range($paraA.Range.Start, $paraA.Range.End - 1)
$rngA.InsertXML((New-PkgXml('<w:p><w:r><w:rPr>' + $rFonts + '</w:rPr><w:t>$ git status</w:t></w:r></w:p>')))

# Paragraph B: insert the "$ git add ." run plus the relocated _GoBack bookmark.
$rngB = $d.Range($paraB.Range.Start, $paraB.Range.End - 1)
$rngB.InsertXML((New-PkgXml('<w:p><w:r><w:rPr>' + $rFonts + '</w:rPr><w:t>$ git add .</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')))

# Paragraph C: replace its whole range (including its own paragraph mark) with two
# new empty Plain Text paragraphs; the engine keeps C itself as a trailing, now
# completely empty, paragraph after the two new ones - exactly matching the diff.
$emptyPlainTextPara = '<w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr>' + $rFonts + '</w:rPr></w:pPr></w:p>'
$bodyC = $emptyPlainTextPara + $emptyPlainTextPara
$paraC.Range.InsertXML((New-PkgXml($bodyC)))

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
